$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 111943997
$ws.Range("B5").Value = 96326
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 219798
$ws.Range("F5").Value = "Skogsknipprot"
$ws.Range("G5").Value = "Epipactis helleborine"
$ws.Range("H5").Value = "(L.) Crantz"
$ws.Range("Q5").Value = 682780.8405377725
$ws.Range("R5").Value = 6694488.393080305

# Row 6
$ws.Range("A6").Value = 111943984
$ws.Range("B6").Value = 99413
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 221235
$ws.Range("F6").Value = "Vårärt"
$ws.Range("G6").Value = "Lathyrus vernus"
$ws.Range("H6").Value = "(L.) Bernh."
$ws.Range("Q6").Value = 682929.3627028114
$ws.Range("R6").Value = 6694685.271877083

# Row 7
$ws.Range("A7").Value = 111943980
$ws.Range("B7").Value = 89183
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 3215
$ws.Range("F7").Value = "Rödgul trumpetsvamp"
$ws.Range("G7").Value = "Craterellus lutescens"
$ws.Range("H7").Value = "(Fr.) Fr."
$ws.Range("Q7").Value = 682877.1417635784
$ws.Range("R7").Value = 6694410.432217407

# Row 8
$ws.Range("A8").Value = 111943988
$ws.Range("B8").Value = 107033
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 220320
$ws.Range("F8").Value = "Ängsskära"
$ws.Range("G8").Value = "Serratula tinctoria"
$ws.Range("H8").Value = "L."
$ws.Range("Q8").Value = 682930.0967543643
$ws.Range("R8").Value = 6694720.015570021

# Row 9
$ws.Range("A9").Value = 111943998
$ws.Range("B9").Value = 98535
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 222498
$ws.Range("F9").Value = "Blåsippa"
$ws.Range("G9").Value = "Hepatica nobilis"
$ws.Range("H9").Value = "Schreb."
$ws.Range("Q9").Value = 682757.1772001419
$ws.Range("R9").Value = 6694405.884787144

# Row 10
$ws.Range("A10").Value = 111943990
$ws.Range("B10").Value = 101703
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 222412
$ws.Range("F10").Value = "Tibast"
$ws.Range("G10").Value = "Daphne mezereum"
$ws.Range("H10").Value = "L."
$ws.Range("Q10").Value = 682930.0967543643
$ws.Range("R10").Value = 6694720.015570021

# Row 11
$ws.Range("A11").Value = 111943983
$ws.Range("B11").Value = 90678
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 4366
$ws.Range("F11").Value = "Skarp dropptaggsvamp"
$ws.Range("G11").Value = "Hydnellum peckii"
$ws.Range("H11").Value = "Banker"
$ws.Range("Q11").Value = 682871.1304590552
$ws.Range("R11").Value = 6694480.539619928

# Row 12
$ws.Range("A12").Value = 111943992
$ws.Range("B12").Value = 89183
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 3215
$ws.Range("F12").Value = "Rödgul trumpetsvamp"
$ws.Range("G12").Value = "Craterellus lutescens"
$ws.Range("H12").Value = "(Fr.) Fr."
$ws.Range("Q12").Value = 682866.8554180798
$ws.Range("R12").Value = 6694644.443727687

# Row 13
$ws.Range("A13").Value = 111943979
$ws.Range("B13").Value = 96253
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 504
$ws.Range("F13").Value = "Guckusko"
$ws.Range("G13").Value = "Cypripedium calceolus"
$ws.Range("H13").Value = "L."
$ws.Range("Q13").Value = 682878.8271195606
$ws.Range("R13").Value = 6694406.550233844

# Row 14
$ws.Range("A14").Value = 111943999
$ws.Range("B14").Value = 99413
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 221235
$ws.Range("F14").Value = "Vårärt"
$ws.Range("G14").Value = "Lathyrus vernus"
$ws.Range("H14").Value = "(L.) Bernh."
$ws.Range("Q14").Value = 682757.1772001419
$ws.Range("R14").Value = 6694405.884787144

# Row 15
$ws.Range("A15").Value = 111943995
$ws.Range("B15").Value = 88899
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 3286
$ws.Range("F15").Value = "Flattoppad klubbsvamp"
$ws.Range("G15").Value = "Clavariadelphus truncatus"
$ws.Range("H15").Value = "(Quél.) Donk"
$ws.Range("Q15").Value = 682779.1674098044
$ws.Range("R15").Value = 6694551.279700429

# Row 16
$ws.Range("A16").Value = 111943981
$ws.Range("B16").Value = 96253
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 504
$ws.Range("F16").Value = "Guckusko"
$ws.Range("G16").Value = "Cypripedium calceolus"
$ws.Range("H16").Value = "L."
$ws.Range("Q16").Value = 682877.1417635784
$ws.Range("R16").Value = 6694410.432217407

# Row 17
$ws.Range("A17").Value = 111943996
$ws.Range("B17").Value = 90332
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 4769
$ws.Range("F17").Value = "Svavelriska"
$ws.Range("G17").Value = "Lactarius scrobiculatus"
$ws.Range("H17").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q17").Value = 682785.3360249697
$ws.Range("R17").Value = 6694547.127516991

# Move the "Alder-Stadium" (K) annotation from row 14 to row 5
$ws.Range("K14").ClearContents()
$ws.Range("K5").Value = "i frukt"
